# Applies the "suppression page etude complementaire" change:
#  - Updates the Metadata "Date" value
#  - Replaces the "Operation / constraint / =" SNOMED filter row on the
#    "Include #0" sheet with a simple "Codes / All codes" row
#  - Updates the system URI from SNOMED CT to EDQM Standard Terms

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the Date property ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-04-09T16:24:06+00:00"

# --- Include #0 sheet: drop the filter columns, switch to "all codes" ---
$wsInclude = $wb.Worksheets.Item("Include #0")

# Fully remove the now-unused B/C cells on rows 1-2 (Operation/Value columns
# that described the SNOMED CT "constraint" filter) without shifting the
# rows below - Clear() drops the <c> nodes entirely, unlike ClearContents().
$wsInclude.Range("B1:C2").Clear()

$wsInclude.Range("A1").Value = "Codes"
$wsInclude.Range("A2").Value = "All codes"

# Swap the referenced code system from SNOMED CT to EDQM Standard Terms.
$wsInclude.Range("B4").Value = "http://standardterms.edqm.eu"
